# Update sheet #5 ("債務" / debt sheet): add species/debtor columns to the
# sharedStrings table, add new property_category/category/date/legislator_name/
# legislator_id/source_file/index columns (H:N) and populate row 2 & 3 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"

# Give the freshly-created header cells (H1:N1) the same formatting
# (bold, centered, bordered) already used by B1:G1.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Row 2 ----
$ws.Range("B2").Value = "現金"
$ws.Range("C2").Value = "高金素梅"
$ws.Range("D2").Value = "陳麗卿新北市泰山區明志路"
$ws.Range("E2").Value = 6000000
$ws.Range("F2").Value = "96年02月06日"
$ws.Range("G2").Value = "借款"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
# "date" column holds a plain ISO-style text string -- force text format so
# Excel doesn't reinterpret it as a date serial number, then restore the
# regular (unbordered) number formatting used by the rest of the row.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-22"
$ws.Range("F2").Copy()
$ws.Range("J2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K2").Value = "高金素梅"
$ws.Range("L2").Value = 926
$ws.Range("M2").Value = "tmp2f3b1"
$ws.Range("N2").Value = 84

# ---- Row 3 ----
$ws.Range("B3").Value = "現金"
$ws.Range("C3").Value = "高金素梅"
$ws.Range("D3").Value = "石旭松新北市泰山區明志路"
$ws.Range("E3").Value = 4000000
$ws.Range("F3").Value = "96年02月06日"
$ws.Range("G3").Value = "借款"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2011-11-22"
$ws.Range("F3").Copy()
$ws.Range("J3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("K3").Value = "高金素梅"
$ws.Range("L3").Value = 926
$ws.Range("M3").Value = "tmp2f3b1"
$ws.Range("N3").Value = 85
